# "removed demo script from m2 & added to m1"
#
# This deck (Module 2 - Setting up the Environments) contains three
# "demo" placeholder slides that introduce a live demo segment:
#   - Slide 10: "Obtaining an Office 365 Developer Subscription" / demo
#   - Slide 15: "Obtaining a Windows Azure trial Subscription" / demo
#   - Slide 31: "Creating a Provider-Hosted App" / demo
#
# The author pulled the demo script out of this module (m2) and moved
# it into module 1. Here that shows up as those three slides being
# hidden from the slide show (they stay in the deck, but are marked
# "don't show" so they're skipped when presenting).

$p = $ppt.ActivePresentation

$demoSlideIndexes = @(10, 15, 31)

foreach ($idx in $demoSlideIndexes) {
    $slide = $p.Slides.Item($idx)
    $slide.SlideShowTransition.Hidden = $true
}
